$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Remove the three rows (old rows 15-17, values for weeks 2023-08-06, 08-13, 08-20)
# which shifts rows 18-60 up to become rows 15-57.
$ws1.Rows("15:17").Delete()

# Row 14's requested-quantity value changes 16 -> 8
$ws1.Range("B14").Value = 8

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B6").Value = 8
